$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D, J, K, L, M, P on rows 2-9, following the
# weekly re-ordering described in the commit ("Fruta / hortaliza, semanal").
$data = @{
    2 = @{ D = 45084; J = 90;  K = 22000; L = 23000; M = 22556; P = 1504 }
    3 = @{ D = 45119; J = 50;  K = 20000; L = 20000; M = 20000; P = 1333 }
    4 = @{ D = 45141; J = 50;  K = 8500;  L = 9000;  M = 8800;  P = 587  }
    5 = @{ D = 44750; J = 140; K = 19000; L = 20000; M = 19571; P = 1305 }
    6 = @{ D = 44749; J = 90;  K = 17000; L = 18000; M = 17556; P = 1170 }
    7 = @{ D = 45133; J = 50;  K = 22000; L = 22000; M = 22000; P = 1467 }
    8 = @{ D = 44839; J = 50;  K = 15000; L = 16000; M = 15600; P = 1040 }
    9 = @{ D = 45091; J = 40;  K = 20000; L = 22000; M = 21000; P = 1400 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
